# Append new scrape results (2026-01-24 06:28 JST) to the "ランサーズ" sheet.
# Rows 2-5 are overwritten with the four newest listings, the previously
# newest-4 rows' data shifts out (rows 6-16 are dropped entirely), and the
# two data-driven column widths (D, H) are narrowed slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Overwrite rows 2-5 with the new scrape snapshot ------------------

$ws.Range("A2").Value = "2026-01-24 06:28:54"
$ws.Range("B2").Value = "【AI×業務改善】AI業務改善のプロ募集|設計サポート+実装(バイブコーディング)【月額固定+時給】"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5477958"
$ws.Range("G2").Value = 338
$ws.Range("H2").Value = "🔥AI,Ai ◇業務改善"

$ws.Range("A3").Value = "2026-01-24 06:28:54"
$ws.Range("B3").Value = "【急募】Amazon SP-API 自動化開発者を探しています"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5477903"
$ws.Range("G3").Value = 328
$ws.Range("H3").Value = "🔥API ◆開発,自動化"

$ws.Range("A4").Value = "2026-01-24 06:28:54"
$ws.Range("B4").Value = "【急募】新しいJob matching appの開発を依頼したい"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5477985"
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = "◆開発"

$ws.Range("A5").Value = "2026-01-24 06:28:54"
$ws.Range("B5").Value = "進行管理およびチームディレクションを担当"
$ws.Range("D5").Value = "~ 5,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = "◇管理"

# --- 2. Point the F2:F5 hyperlinks at the matching new URLs --------------
# (Collect into a plain array first: deleting/updating while iterating the
# live Hyperlinks collection shifts indices out from under a foreach.)

$links = @()
foreach ($h in $ws.Hyperlinks) {
    $links += $h
}

$links[0].Address = "https://www.lancers.jp/work/detail/5477958"
$links[1].Address = "https://www.lancers.jp/work/detail/5477903"
$links[2].Address = "https://www.lancers.jp/work/detail/5477985"
$links[3].Address = "https://www.lancers.jp/work/detail/5418064"

# --- 3. Drop the old rows 6-16 (11 rows) so the sheet ends at row 5 -------

$ws.Rows.Item(6).Resize(11).Delete()

# --- 4. Remove the now-orphaned hyperlinks that used to target F6:F16 ----

$remaining = @()
foreach ($h in $ws.Hyperlinks) {
    $remaining += $h
}
for ($i = $remaining.Count - 1; $i -ge 4; $i--) {
    $remaining[$i].Delete()
}

# --- 5. Narrow columns D (32 -> 30) and H (19 -> 14) ----------------------
# ColumnWidth uses Excel's font-metric units, which differ from the raw
# stored <col width>; compute the engine's current offset from the
# before-edit stored widths so the saved XML lands on exactly 30 / 14.

$curD = $ws.Columns.Item(4).ColumnWidth()
$curH = $ws.Columns.Item(8).ColumnWidth()
$offsetD = 32 - $curD
$offsetH = 19 - $curH
$ws.Columns.Item(4).ColumnWidth = 30 - $offsetD
$ws.Columns.Item(8).ColumnWidth = 14 - $offsetH
